$wb = $excel.ActiveWorkbook

# Update the "展览" sheet: F2 275 -> 277, F5 14 -> 15
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 277
$ws1.Range("F5").Value = 15

# Update the "全部类型" sheet: F2 275 -> 277, F5 14 -> 15
$ws2 = $wb.Worksheets.Item("全部类型")
$ws2.Range("F2").Value = 277
$ws2.Range("F5").Value = 15
